$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 19.4296083400938
$ws.Range("D2").Value = 8.374397451977321
$ws.Range("E2").Value = 25.48132220474869
$ws.Range("F2").Value = 47.30507581188667
$ws.Range("G2").Value = 3.657740745787851
$ws.Range("L2").Value = 11.66502544331799
$ws.Range("M2").Value = 18.04384958489591
$ws.Range("N2").Value = 20.07103629680555
$ws.Range("B3").Value = 19.20537247891636
$ws.Range("D3").Value = 8.283703690469107
$ws.Range("E3").Value = 23.97780679880245
$ws.Range("F3").Value = 46.08514827313272
$ws.Range("G3").Value = 3.666295870156662
$ws.Range("L3").Value = 11.49733031449543
$ws.Range("M3").Value = 17.9172231198017
$ws.Range("N3").Value = 20.16484919806295
$ws.Range("B4").Value = 19.07202451056637
$ws.Range("D4").Value = 8.231340986535637
$ws.Range("E4").Value = 23.00453649840849
$ws.Range("F4").Value = 45.34269954775912
$ws.Range("G4").Value = 3.671794546835565
$ws.Range("L4").Value = 11.39683617571893
$ws.Range("M4").Value = 17.84390406800298
$ws.Range("N4").Value = 20.22468865941824
$ws.Range("B5").Value = 19.0188293588224
$ws.Range("D5").Value = 8.210852331281316
$ws.Range("E5").Value = 22.59540845639168
$ws.Range("F5").Value = 45.04227397573622
$ws.Range("G5").Value = 3.674097539858077
$ws.Range("L5").Value = 11.35654994593969
$ws.Range("M5").Value = 17.81516040125308
$ws.Range("N5").Value = 20.24964054867206
$ws.Range("B6").Value = 19.0100670847468
$ws.Range("D6").Value = 8.207501875948118
$ws.Range("E6").Value = 22.52672045923709
$ws.Range("F6").Value = 44.99253098545147
$ws.Range("G6").Value = 3.674483721619219
$ws.Range("L6").Value = 11.34990194169884
$ws.Range("M6").Value = 17.81045661326005
$ws.Range("N6").Value = 20.25381814398938
$ws.Range("B7").Value = 19.07130239370382
$ws.Range("D7").Value = 8.23106121291911
$ws.Range("E7").Value = 22.99906938098054
$ws.Range("F7").Value = 45.33863865817002
$ws.Range("G7").Value = 3.671825353237795
$ws.Range("L7").Value = 11.39629010738644
$ws.Range("M7").Value = 17.84351180197579
$ws.Range("N7").Value = 20.22502286889866
$ws.Range("B8").Value = 19.35142635497798
$ws.Range("D8").Value = 8.342445403496946
$ws.Range("E8").Value = 24.97335296492073
$ws.Range("F8").Value = 46.88335624390892
$ws.Range("G8").Value = 3.660639798661002
$ws.Range("L8").Value = 11.60671834332674
$ws.Range("M8").Value = 17.99928217366265
$ws.Range("N8").Value = 20.10292096847113
$ws.Range("B9").Value = 19.9325911361997
$ws.Range("D9").Value = 8.586528661354798
$ws.Range("E9").Value = 28.44490083935823
$ws.Range("F9").Value = 49.94476157392702
$ws.Range("G9").Value = 3.640635450799039
$ws.Range("L9").Value = 12.03696584035622
$ws.Range("M9").Value = 18.33891149995984
$ws.Range("N9").Value = 19.88104769582826
$ws.Range("B10").Value = 20.37546712568976
$ws.Range("D10").Value = 8.780415991974454
$ws.Range("E10").Value = 30.74986890657146
$ws.Range("F10").Value = 52.18799880056635
$ws.Range("G10").Value = 3.627087166531943
$ws.Range("L10").Value = 12.36106402615421
$ws.Range("M10").Value = 18.60782969376585
$ws.Range("N10").Value = 19.72847759161008
$ws.Range("B11").Value = 20.57962780311637
$ws.Range("D11").Value = 8.871519672075125
$ws.Range("E11").Value = 31.745111688737
$ws.Range("F11").Value = 53.20224901820757
$ws.Range("G11").Value = 3.621167037274998
$ws.Range("L11").Value = 12.50966034939157
$ws.Range("M11").Value = 18.73404274221526
$ws.Range("N11").Value = 19.66127710290276
$ws.Range("B12").Value = 20.65726108996916
$ws.Range("D12").Value = 8.906412100436979
$ws.Range("E12").Value = 32.11433174818437
$ws.Range("F12").Value = 53.58504999487188
$ws.Range("G12").Value = 3.618959693189994
$ws.Range("L12").Value = 12.56604931999865
$ws.Range("M12").Value = 18.78236432112084
$ws.Range("N12").Value = 19.63614227875722
$ws.Range("B13").Value = 20.64052808269719
$ws.Range("D13").Value = 8.898880283435567
$ws.Range("E13").Value = 32.03515391463768
$ws.Range("F13").Value = 53.50266933934449
$ws.Range("G13").Value = 3.619433558285649
$ws.Range("L13").Value = 12.55390042174855
$ws.Range("M13").Value = 18.77193444314624
$ws.Range("N13").Value = 19.64154168413844
$ws.Range("B14").Value = 20.58600865790637
$ws.Range("D14").Value = 8.874382548185238
$ws.Range("E14").Value = 31.77564101988398
$ws.Range("F14").Value = 53.23376988291749
$ws.Range("G14").Value = 3.620984749436335
$ws.Range("L14").Value = 12.51429739663465
$ws.Range("M14").Value = 18.73800779543551
$ws.Range("N14").Value = 19.65920300777073
$ws.Range("B15").Value = 20.55265394890376
$ws.Range("D15").Value = 8.859427489926848
$ws.Range("E15").Value = 31.61568513384571
$ws.Range("F15").Value = 53.06888455916234
$ws.Range("G15").Value = 3.62193937506791
$ws.Range("L15").Value = 12.4900534511543
$ws.Range("M15").Value = 18.7172945084051
$ws.Range("N15").Value = 19.6700616472859
$ws.Range("B16").Value = 20.3621733614043
$ws.Range("D16").Value = 8.774518535571163
$ws.Range("E16").Value = 30.68375535140904
$ws.Range("F16").Value = 52.1215569933217
$ws.Range("G16").Value = 3.627478911725154
$ws.Range("L16").Value = 12.35137213612204
$ws.Range("M16").Value = 18.5996568690836
$ws.Range("N16").Value = 19.73291328143227
$ws.Range("B17").Value = 20.24596309762819
$ws.Range("D17").Value = 8.723156390938067
$ws.Range("E17").Value = 30.09839756027127
$ws.Range("F17").Value = 51.5385315646076
$ws.Range("G17").Value = 3.630939157602547
$ws.Range("L17").Value = 12.2665585973611
$ws.Range("M17").Value = 18.52846246371796
$ws.Range("N17").Value = 19.77203218488492
$ws.Range("B18").Value = 20.17937951159827
$ws.Range("D18").Value = 8.693889020037334
$ws.Range("E18").Value = 29.75669733940934
$ws.Range("F18").Value = 51.2026272017354
$ws.Range("G18").Value = 3.632952305432883
$ws.Range("L18").Value = 12.21788834390992
$ws.Range("M18").Value = 18.48788050277523
$ws.Range("N18").Value = 19.79474010431233
$ws.Range("B19").Value = 20.15688167528618
$ws.Range("D19").Value = 8.684027504836525
$ws.Range("E19").Value = 29.64014162690026
$ws.Range("F19").Value = 51.08881104858759
$ws.Range("G19").Value = 3.633637871012703
$ws.Range("L19").Value = 12.20143020392536
$ws.Range("M19").Value = 18.47420412932519
$ws.Range("N19").Value = 19.80246443311527
$ws.Range("B20").Value = 20.25830770498326
$ws.Range("D20").Value = 8.728595712599617
$ws.Range("E20").Value = 30.16122920780915
$ws.Range("F20").Value = 51.60065679471937
$ws.Range("G20").Value = 3.630568441228901
$ws.Range("L20").Value = 12.27557587995141
$ws.Range("M20").Value = 18.53600344717775
$ws.Range("N20").Value = 19.76784644037353
$ws.Range("B21").Value = 20.60201411157177
$ws.Range("D21").Value = 8.88156764408817
$ws.Range("E21").Value = 31.85207394158859
$ws.Range("F21").Value = 53.3127896206931
$ws.Range("G21").Value = 3.620528194993109
$ws.Range("L21").Value = 12.52592692151062
$ws.Range("M21").Value = 18.74795880799172
$ws.Range("N21").Value = 19.65400700241641
$ws.Range("B22").Value = 20.8284935758129
$ws.Range("D22").Value = 8.983825350790713
$ws.Range("E22").Value = 32.91251930810464
$ws.Range("F22").Value = 54.42418884131107
$ws.Range("G22").Value = 3.614167067593289
$ws.Range("L22").Value = 12.69021653308949
$ws.Range("M22").Value = 18.88954343922444
$ws.Range("N22").Value = 19.5814259645938
$ws.Range("B23").Value = 20.7074695989719
$ws.Range("D23").Value = 8.929047821214823
$ws.Range("E23").Value = 32.35061735490881
$ws.Range("F23").Value = 53.83182556619604
$ws.Range("G23").Value = 3.617543909533934
$ws.Range("L23").Value = 12.60248636718769
$ws.Range("M23").Value = 18.81370743231908
$ws.Range("N23").Value = 19.61999882606781
$ws.Range("B24").Value = 20.25272599522516
$ws.Range("D24").Value = 8.726135779224322
$ws.Range("E24").Value = 30.1328391245293
$ws.Range("F24").Value = 51.57257217689821
$ws.Range("G24").Value = 3.630735967914227
$ws.Range("L24").Value = 12.27149888191836
$ws.Range("M24").Value = 18.53259308027322
$ws.Range("N24").Value = 19.76973813629468
$ws.Range("B25").Value = 19.7723203567181
$ws.Range("D25").Value = 8.517855902493238
$ws.Range("E25").Value = 27.54882624627342
$ws.Range("F25").Value = 49.11584911121955
$ws.Range("G25").Value = 3.64584335243088
$ws.Range("L25").Value = 11.91896979571501
$ws.Range("M25").Value = 18.24351719036726
$ws.Range("N25").Value = 19.93921716925843
